# Consolidate input metafile file contents to cover uploader and rcsb importer.
# Splits the "PDB filename" column into a base ID column and a file-extension
# column, relabels the "Is model" column, and appends a new "From RCSB" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$ws.Range("D1").Value = "PDB or RCSB ID"
$ws.Range("E1").Value = "File Extension"
$ws.Range("F1").Value = " Is model"
$ws.Range("G1").Value = "From RCSB"

# --- Data rows ----------------------------------------------------------
# row -> base id (without .pdb extension)
$baseIds = @{
    2 = "JCVISYN3_0001"
    3 = "JCVISYN3_0001"
    4 = "JCVISYN3_0003"
    5 = "JCVISYN3_0004"
}

for ($r = 2; $r -le 5; $r++) {
    $isModel = $ws.Cells.Item($r, 5).Text
    $ws.Cells.Item($r, 4).Value = $baseIds[$r]
    $ws.Cells.Item($r, 5).Value = "pdb"
    $ws.Cells.Item($r, 6).Value = $isModel
    $ws.Cells.Item($r, 7).Value = "n"
}

# --- Column widths --------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 21

# --- Selection / active cell --------------------------------------------
$ws.Range("G6").Select()
